$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh crypto list: updated prices and 1h volume-change percentages.
# Rows 18-20 are also reordered (BitcoinCash, ShibaInu, Chainlink) to reflect
# the coins' new ranking order.
# A leading apostrophe is used for purely-numeric-looking price strings so
# Excel keeps them as text (matching the source data format) instead of
# auto-converting them to numbers.
$ws.Range("D2").Value = "27.459.22"
$ws.Range("E2").Value = "  -0.77%  "
$ws.Range("D3").Value = "1.573.36"
$ws.Range("E3").Value = "  -0.65%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").Value = "'207.58"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").Value = "'0.498"
$ws.Range("E6").Value = "  -0.81%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("D10").Value = "'0.0595"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("D12").Value = "1.797.63"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").Value = "1.573.34"
$ws.Range("E13").Value = "  -0.74%  "
$ws.Range("E14").Value = "  -1.15%  "
$ws.Range("E15").Value = "  -1.79%  "
$ws.Range("D16").Value = "'63.50"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("D17").Value = "27.459.90"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("B18").Value = "BitcoinCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D18").Value = "'214.32"
$ws.Range("E18").Value = "  -2.40%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "0.0₃0693"
$ws.Range("E19").Value = "  -0.10%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'7.31"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("E21").Value = "  -0.29%  "
$ws.Range("D22").Value = "'4.14"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").Value = "'9.49"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("E24").Value = "  +1.10%  "
$ws.Range("D25").Value = "'153.20"
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "'6.68"
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").Value = "'14.95"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("E28").Value = "  -0.74%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("E32").Value = "  -0.60%  "
$ws.Range("D33").Value = "1.395.15"
$ws.Range("E33").Value = "  +2.38%  "
$ws.Range("E34").Value = "  +0.95%  "
$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  +2.39%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.29"
$ws.Range("E36").Value = "  -0.86%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'0.941"
$ws.Range("E37").Value = "  -3.44%  "
$ws.Range("E38").Value = "  -0.04%  "
$ws.Range("D39").Value = "'0.532"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "'0.824"
$ws.Range("E40").Value = "  +0.62%  "
$ws.Range("E41").Value = "  -0.26%  "
$ws.Range("D42").Value = "'0.995"
$ws.Range("E42").Value = "  +1.49%  "
$ws.Range("D43").Value = "'1.82"
$ws.Range("E43").Value = "  +5.03%  "
$ws.Range("D44").Value = "'64.45"
$ws.Range("E44").Value = "  +1.31%  "
$ws.Range("D45").Value = "'2.18"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("D46").Value = "'5.25"
$ws.Range("E46").Value = "  +1.40%  "
$ws.Range("D47").Value = "1.709.10"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").Value = "'85.95"
$ws.Range("E48").Value = "  -2.27%  "
$ws.Range("D49").Value = "0.0₆0100"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").Value = "'0.0955"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("E51").Value = "  -0.47%  "
